# Update "want to go" (F column) counts across sheets, per commit
# "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1799
$ws1.Range("F8").Value  = 345
$ws1.Range("F9").Value  = 1754
$ws1.Range("F12").Value = 818
$ws1.Range("F13").Value = 344
$ws1.Range("F14").Value = 689
$ws1.Range("F15").Value = 12867
$ws1.Range("F16").Value = 12849
$ws1.Range("F22").Value = 583
$ws1.Range("F27").Value = 79

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 54

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1799
$ws4.Range("F12").Value = 54
$ws4.Range("F13").Value = 345
$ws4.Range("F14").Value = 1754
$ws4.Range("F17").Value = 818
$ws4.Range("F18").Value = 344
$ws4.Range("F20").Value = 689
$ws4.Range("F21").Value = 12867
$ws4.Range("F22").Value = 12849
$ws4.Range("F28").Value = 583
$ws4.Range("F37").Value = 79

$wb.Save()
